# Updated symbol list on Sat Feb 11 22:32:18 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) values for
# the coin rows on Sheet1. The source data stores these as plain text
# (e.g. "309.81", "0.96%") rather than numbers, so each target cell is
# first formatted as Text ("@") before its value is assigned. This keeps
# Excel from re-interpreting the numeric-looking / percent-looking
# strings as actual numbers (which would silently change their stored
# representation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ D = "<new price>" (optional); E = "<new volume%>" }
$updates = [ordered]@{
    2  = @{ D = "310.27";        E = "1.22%" }
    3  = @{ D = "41.04";         E = "1.88%" }
    4  = @{ D = "5.121";         E = "0.47%" }
    5  = @{ D = "0.07680";       E = "1.30%" }
    6  = @{ D = "4.273";         E = "0.00%" }
    7  = @{ D = "1.620";         E = "1.09%" }
    8  = @{ D = "0.9192";        E = "1.52%" }
    10 = @{ D = "0.1228";        E = "21.40%" }
    11 = @{ D = "0.1827";        E = "4.19%" }
    12 = @{ D = "0.09105";       E = "-0.41%" }
    13 = @{ D = "0.04285";       E = "2.50%" }
    14 = @{ D = "0.1050";        E = "-0.37%" }
    15 = @{ D = "0.001251";      E = "0.10%" }
    16 = @{ D = "0.005634";      E = "-3.35%" }
    17 = @{               E = "0.07%" }
    18 = @{               E = "1.23%" }
    19 = @{ D = "6.921";         E = "3.90%" }
    20 = @{ D = "0.1388";        E = "2.28%" }
    21 = @{ D = "0.2730";        E = "0.05%" }
    22 = @{ D = "0.04043";       E = "-3.37%" }
    23 = @{ D = "0.001266";      E = "3.02%" }
    24 = @{ D = "0.004075";      E = "0.33%" }
    25 = @{ D = "0.0001268";     E = "-2.80%" }
    26 = @{               E = "24.50%" }
    38 = @{ D = "0.02461";       E = "2.20%" }
    39 = @{ D = "0.05255";       E = "2.24%" }
    40 = @{ D = "0.007832";      E = "0.57%" }
    41 = @{ D = "0.1312";        E = "1.60%" }
    42 = @{ D = "0.006794";      E = "-3.73%" }
    43 = @{ D = "0.001839";      E = "-5.42%" }
    44 = @{ D = "0.008174";      E = "-3.46%" }
    45 = @{ D = "0.3097";        E = "-6.79%" }
    46 = @{ D = "0.00006832";    E = "7.33%" }
    47 = @{ D = "0.00000000749"; E = "-0.51%" }
    48 = @{ D = "0.2232";        E = "604.84%" }
    49 = @{ D = "0.004092";      E = "-7.17%" }
    50 = @{ D = "0.00002096";    E = "-0.51%" }
    51 = @{               E = "-0.51%" }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    foreach ($col in @("D", "E")) {
        if ($vals.Contains($col)) {
            $addr = "$col$row"
            $ws.Range($addr).NumberFormat = "@"
            $ws.Range($addr).Value = $vals[$col]
        }
    }
}
